# Insert a new row at position 34, shifting existing rows 34-112 down to 35-113,
# then populate the newly inserted row 34 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(34).Insert()

$ws.Cells.Item(34, 1).Value = 5
$ws.Cells.Item(34, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(34, 3).Value = "Maule"
$ws.Cells.Item(34, 4).Value = [DateTime]"2022-11-18"
$ws.Cells.Item(34, 5).Value = 7
$ws.Cells.Item(34, 6).Value = 100112026
$ws.Cells.Item(34, 7).Value = "Haba"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 500
$ws.Cells.Item(34, 11).Value = 8000
$ws.Cells.Item(34, 12).Value = 8000
$ws.Cells.Item(34, 13).Value = 8000
$ws.Cells.Item(34, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(34, 15).Value = "Región del Maule"
$ws.Cells.Item(34, 16).Value = 320
$ws.Cells.Item(34, 17).Value = 25
$ws.Cells.Item(34, 18).Value = "Hortaliza"
